{"js": "// The document body contains 101 paragraphs in document order:\n// index 0 is the title date line, indices 1..100 are the 100 table-cell\n// equations (5 columns x 20 rows, row-major). This mirrors the edit in\n// the diff, which rewrites each <w:t> text in place without touching\n// any other paragraph/run formatting.\nconst newValues = [\n  \"2024-03-14 Thursday\",\n  \"42-3=39\",\n  \"49+8=57\",\n  \"4-3=1\",\n  \"32+34=66\",\n  \"64-39=25\",\n  \"49+6=55\",\n  \"28-25=3\",\n  \"82-72=10\",\n  \"99-12=87\",\n  \"28+33=61\",\n  \"64-53=11\",\n  \"85-9=76\",\n  \"56-54=2\",\n  \"92-55=37\",\n  \"70+0=70\",\n  \"37+0=37\",\n  \"51-8=43\",\n  \"85-34=51\",\n  \"42+33=75\",\n  \"15+35=50\",\n  \"71-63=8\",\n  \"18+54=72\",\n  \"97-88=9\",\n  \"97-9=88\",\n  \"68-0=68\",\n  \"82-51=31\",\n  \"40+17=57\",\n  \"24+50=74\",\n  \"5+61=66\",\n  \"18+10=28\",\n  \"32+5=37\",\n  \"99-30=69\",\n  \"12-4=8\",\n  \"28-6=22\",\n  \"38+26=64\",\n  \"79-13=66\",\n  \"5+74=79\",\n  \"43-30=13\",\n  \"60-57=3\",\n  \"87+3=90\",\n  \"41-8=33\",\n  \"92-27=65\",\n  \"91-79=12\",\n  \"25+14=39\",\n  \"62+35=97\",\n  \"50-2=48\",\n  \"56+21=77\",\n  \"92+4=96\",\n  \"7+88=95\",\n  \"42-8=34\",\n  \"68-29=39\",\n  \"15+18=33\",\n  \"59-5=54\",\n  \"35+4=39\",\n  \"20+4=24\",\n  \"0+18=18\",\n  \"5+76=81\",\n  \"91+4=95\",\n  \"35-10=25\",\n  \"36+24=60\",\n  \"12+58=70\",\n  \"21+23=44\",\n  \"5+76=81\",\n  \"87-71=16\",\n  \"13+6=19\",\n  \"92-12=80\",\n  \"82-0=82\",\n  \"0+95=95\",\n  \"78+2=80\",\n  \"79-76=3\",\n  \"41-19=22\",\n  \"77+20=97\",\n  \"59-3=56\",\n  \"43-9=34\",\n  \"90-58=32\",\n  \"67+2=69\",\n  \"55+3=58\",\n  \"14+7=21\",\n  \"76-73=3\",\n  \"74-52=22\",\n  \"97-59=38\",\n  \"65+6=71\",\n  \"64+3=67\",\n  \"22-16=6\",\n  \"57-16=41\",\n  \"45-43=2\",\n  \"62+7=69\",\n  \"23-0=23\",\n  \"22+58=80\",\n  \"66+16=82\",\n  \"2+31=33\",\n  \"53-53=0\",\n  \"67-60=7\",\n  \"52-21=31\",\n  \"62-20=42\",\n  \"17+65=82\",\n  \"18+32=50\",\n  \"10+79=89\",\n  \"81-21=60\",\n  \"69-1=68\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + paragraphs.items.length +\n    \" (expected \" + newValues.length + \")\"\n  );\n}\n\nfor (let i = 0; i < newValues.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document has a title paragraph (the date line) followed by a single\n# 20-row x 5-column table where every cell holds one equation. This mirrors\n# the diff exactly: only the <w:t> text content changes, in document order,\n# run/paragraph formatting is left untouched.\n\n$d = $word.ActiveDocument\n\n$titleText = \"2024-03-14 Thursday\"\n$d.Paragraphs.Item(1).Range.Text = $titleText\n\n$newValues = @(\n    \"42-3=39\", \"49+8=57\", \"4-3=1\", \"32+34=66\", \"64-39=25\",\n    \"49+6=55\", \"28-25=3\", \"82-72=10\", \"99-12=87\", \"28+33=61\",\n    \"64-53=11\", \"85-9=76\", \"56-54=2\", \"92-55=37\", \"70+0=70\",\n    \"37+0=37\", \"51-8=43\", \"85-34=51\", \"42+33=75\", \"15+35=50\",\n    \"71-63=8\", \"18+54=72\", \"97-88=9\", \"97-9=88\", \"68-0=68\",\n    \"82-51=31\", \"40+17=57\", \"24+50=74\", \"5+61=66\", \"18+10=28\",\n    \"32+5=37\", \"99-30=69\", \"12-4=8\", \"28-6=22\", \"38+26=64\",\n    \"79-13=66\", \"5+74=79\", \"43-30=13\", \"60-57=3\", \"87+3=90\",\n    \"41-8=33\", \"92-27=65\", \"91-79=12\", \"25+14=39\", \"62+35=97\",\n    \"50-2=48\", \"56+21=77\", \"92+4=96\", \"7+88=95\", \"42-8=34\",\n    \"68-29=39\", \"15+18=33\", \"59-5=54\", \"35+4=39\", \"20+4=24\",\n    \"0+18=18\", \"5+76=81\", \"91+4=95\", \"35-10=25\", \"36+24=60\",\n    \"12+58=70\", \"21+23=44\", \"5+76=81\", \"87-71=16\", \"13+6=19\",\n    \"92-12=80\", \"82-0=82\", \"0+95=95\", \"78+2=80\", \"79-76=3\",\n    \"41-19=22\", \"77+20=97\", \"59-3=56\", \"43-9=34\", \"90-58=32\",\n    \"67+2=69\", \"55+3=58\", \"14+7=21\", \"76-73=3\", \"74-52=22\",\n    \"97-59=38\", \"65+6=71\", \"64+3=67\", \"22-16=6\", \"57-16=41\",\n    \"45-43=2\", \"62+7=69\", \"23-0=23\", \"22+58=80\", \"66+16=82\",\n    \"2+31=33\", \"53-53=0\", \"67-60=7\", \"52-21=31\", \"62-20=42\",\n    \"17+65=82\", \"18+32=50\", \"10+79=89\", \"81-21=60\", \"69-1=68\"\n)\n\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\nif (($rowCount * $colCount) -ne $newValues.Count) {\n    throw \"Unexpected table size: $rowCount x $colCount (expected $($newValues.Count) cells)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\nWrite-Output \"done\"\n"}
